$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-type the header row text without the old fixed-width trailing padding.
$ws.Range("A1").Value = "Stock Code"
$ws.Range("B1").Value = " Company Name"
$ws.Range("C1").Value = " Sector"
$ws.Range("D1").Value = " Open"
$ws.Range("E1").Value = " Close"
$ws.Range("F1").Value = " Volume"
$ws.Range("G1").Value = " Trade Date"

# The WBC "Open" price (row 4, column D) was corrected from 24.3 to 2.
$ws.Range("D4").Value = 2

$ws.Range("G7").Select() | Out-Null
